$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the standalone "Meta description: ..." paragraph that
#    directly follows the title heading.
# ------------------------------------------------------------------
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2. Replace the closing italic "image prompt" paragraph's text with
#    the (former) meta-description text, and insert a new bold
#    paragraph - duplicating the page title - right before it.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
# Exclude the trailing paragraph mark so InsertXML only replaces the
# run content and adds a fresh paragraph break, instead of leaving a
# stray empty paragraph behind.
$r.MoveEnd(1, -1) | Out-Null

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr><w:b/></w:rPr>
              <w:t>Play Football Cash Collect for Free - Review by Slot Game Writer</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr><w:i/></w:rPr>
              <w:t>Read our review of Football Cash Collect, a Playtech slot game with 5x3 layout, 30 paylines, free spins, and four jackpots. Play for free today.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xmlFrag) | Out-Null
